$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at H: "Amount *" (G) stays, a new "Amount (Folio
#    Currency)" column is inserted immediately after it. Everything that used
#    to live in H (Notes), I (Type) and J (Rule For) shifts one column right
#    (Notes->I, Type->J, Rule For->K). Data validation + widths + cell
#    formatting shift automatically with the insert.
# ---------------------------------------------------------------------------
$ws.Columns("H").Insert()

# ---------------------------------------------------------------------------
# 2. Header text: rename the old "Amount *" header and label the new column.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Amount (Fund Currency)*"
$ws.Range("H1").Value = "Amount (Folio Currency)"

# Match the bold / bordered look used by the other header cells.
$ws.Range("G1:H1").Font.Bold = $true
$ws.Range("G1:H1").Font.Name = "Arial"
$ws.Range("G1:H1").Font.Size = 10
$ws.Range("G1:H1").Borders.Item(9).Weight = -4138

# Widen the two amount columns so the longer headers fit.
$ws.Range("G1").ColumnWidth = 22.6
$ws.Range("H1").ColumnWidth = 22.6

# ---------------------------------------------------------------------------
# 3. Comments are anchored to a fixed cell reference and do NOT shift when a
#    column is inserted, so the comment that used to sit on H1 ("-Optional")
#    and the one that used to sit on I1 ("-Mandatory / nomenclature") are now
#    sitting on the wrong header. Move their text to the cells that now hold
#    the Notes / Type columns (I1 / J1 respectively).
# ---------------------------------------------------------------------------
$optionalText = $ws.Range("H1").Comment.Text()
$nomenclatureText = $ws.Range("I1").Comment.Text()

$ws.Range("H1").Comment.Delete()
$ws.Range("I1").Comment.Delete()

$ws.Range("I1").AddComment($optionalText) | Out-Null
$ws.Range("J1").AddComment($nomenclatureText) | Out-Null

# Re-create the G1 ("Amount (Fund Currency)*") comment as a fresh note from
# the new reviewer so it is attributed to the newly added comment author,
# keeping its original guidance text.
$amountText = $ws.Range("G1").Comment.Text()
$ws.Range("G1").Comment.Delete()
$ws.Range("G1").AddComment($amountText) | Out-Null

# ---------------------------------------------------------------------------
# 4. Misc view bits that Excel records whenever the sheet is touched.
# ---------------------------------------------------------------------------
$ws.Range("H1").Select()
